$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Present-Storage")

$ws.Range("B2").Value = 10.03
$ws.Range("B3").Value = 5.08
$ws.Range("B4").Value = 6.91
$ws.Range("B5").Value = 5.46
$ws.Range("B6").Value = 9.460000000000001
$ws.Range("B7").Value = 4.22
$ws.Range("B8").Value = 6.16
$ws.Range("B9").Value = 4.63
$ws.Range("B10").Value = 17.62
$ws.Range("B11").Value = 13.75
$ws.Range("B12").Value = 15.18
$ws.Range("B13").Value = 14.05
$ws.Range("B14").Value = 10.18
$ws.Range("B15").Value = 5.42
$ws.Range("B16").Value = 7.18
$ws.Range("B17").Value = 5.79
